# fb125/raw/Table4.xlsx - minor changes to some xlsx tables
# (nfishers data almost complete)
#
# 1) The "Los Angeles" row (row 7) was missing its area-of-residence label;
#    fill it in (this also appends a new shared string).
# 2) The active selection was left on D10; move it back to A7, matching
#    where the author was last working.
# 3) Rows 1 and 7 keep their original 16pt height explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Los Angeles"

$ws.Rows.Item(1).RowHeight = 16
$ws.Rows.Item(7).RowHeight = 16

[void]$ws.Range("A7").Select()
